$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 35 / Row 36 coin swap: OKB moves to row 35, PEPE moves to row 36 ---
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "68.43"
$ws.Range("E35").Value = "  -3.07%  "

$ws.Range("B36").Value = "PEPE"
$ws.Range("C36").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0₃0977"
$ws.Range("E36").Value = "  +13.79%  "

# --- Price (column D) updates ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "72.244.83"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.907.22"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.01"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.81"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.678"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.772"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.187"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.83"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.39"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.544.67"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.922.80"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.11"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.03"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "72.360.82"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "443.11"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.74"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "94.92"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.96"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.18"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.12"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.95"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.24"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.44"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.84"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "50.93"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "13.69"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "621.06"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.423"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.23"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0473"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.41"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.90"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.33"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.56"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.822.29"

# --- Volume(1h) percentage (column E) updates ---
$ws.Range("E2").Value = "  -0.69%  "
$ws.Range("E3").Value = "  -2.18%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("E5").Value = "  +1.54%  "
$ws.Range("E6").Value = "  +10.84%  "
$ws.Range("E7").Value = "  -0.71%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  +2.63%  "
$ws.Range("E10").Value = "  +10.56%  "
$ws.Range("E11").Value = "  +1.92%  "
$ws.Range("E12").Value = "  +2.55%  "
$ws.Range("E13").Value = "  +4.44%  "
$ws.Range("E14").Value = "  -1.92%  "
$ws.Range("E15").Value = "  -2.18%  "
$ws.Range("E16").Value = "  +2.61%  "
$ws.Range("E17").Value = "  -0.25%  "
$ws.Range("E19").Value = "  -0.34%  "
$ws.Range("E20").Value = "  -1.27%  "
$ws.Range("E21").Value = "  +1.53%  "
$ws.Range("E22").Value = "  -0.36%  "
$ws.Range("E23").Value = "  -1.55%  "
$ws.Range("E24").Value = "  -5.83%  "
$ws.Range("E25").Value = "  -2.65%  "
$ws.Range("E26").Value = "  -4.44%  "
$ws.Range("E27").Value = "  -3.28%  "
$ws.Range("E28").Value = "  +0.48%  "
$ws.Range("E29").Value = "  -4.24%  "
$ws.Range("E30").Value = "  -2.77%  "
$ws.Range("E31").Value = "  -1.73%  "
$ws.Range("E32").Value = "  +1.59%  "
$ws.Range("E33").Value = "  +0.52%  "
$ws.Range("E34").Value = "  -4.38%  "
$ws.Range("E37").Value = "  -9.00%  "
$ws.Range("E38").Value = "  -4.36%  "
$ws.Range("E39").Value = "  +0.11%  "
$ws.Range("E40").Value = "  +0.19%  "
$ws.Range("E41").Value = "  -2.77%  "
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("E43").Value = "  +42.47%  "
$ws.Range("E44").Value = "  -3.55%  "
$ws.Range("E45").Value = "  -6.42%  "
$ws.Range("E46").Value = "  -2.65%  "
$ws.Range("E47").Value = "  -13.84%  "
$ws.Range("E48").Value = "  -1.51%  "
$ws.Range("E49").Value = "  -9.72%  "
$ws.Range("E50").Value = "  +0.52%  "
$ws.Range("E51").Value = "  +2.40%  "
